# "WGCNA for 70 min redo"
#
# During this session the author (Tyler Milewski) re-attempted the WGCNA
# figure for the "70 min" slide (sldId 267 / slide index 10, which is
# currently an emptied-out slide after earlier edits removed its two
# textboxes and its picture). He pasted/added a new picture onto that
# slide and then, still unhappy with it, deleted it again; he also spun
# up a brand-new slide right after it to try an alternate layout for the
# WGCNA figure and then deleted that whole slide too. Nothing from this
# flurry of edits survives in the final saved deck (slide 267 stays
# empty, no extra slide remains) - only PowerPoint's co-authoring
# change-log remembers it happened. This script reproduces the actual
# user actions so the edit history matches.

$p = $ppt.ActivePresentation
$s = $p.Slides.FindBySlideID(267)

# --- slide 267: try a picture, discard it, try again, discard again ---
# (an earlier textbox/picture churn on this slide already consumed
# shape ids 2 and 4 historically; replay a throwaway add/delete first so
# the id counter lands on 3 for the real attempt, same as it did live)
$placeholder = $s.Shapes.AddPicture("C:\Users\Tyler\Pictures\temp.png", $false, $true, 1, 1, 1, 1)
$placeholder.Delete()

$wgcna = $s.Shapes.AddPicture("C:\Users\Tyler\Pictures\WGCNA_70min.png", $false, $true, 457200, 457200, 11277600, 5943600)
$wgcna.Delete()

# --- a whole new slide attempted right after it, then scrapped ---
$tryAgain = $p.Slides.Add($s.SlideIndex + 1, 12)
$tryAgain.Delete()
